$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44981, 7839385000000),
    @(44988, 7830796000000),
    @(44995, 7829368000000),
    @(45002, 7831859000000),
    @(45009, 7835921000000),
    @(45016, 7729612000000),
    @(45023, 7729251000000),
    @(45030, 7730924000000),
    @(45037, 7714222000000),
    @(45044, 7719816000000),
    @(45051, 7716913000000),
    @(45058, 7728510000000),
    @(45065, 7730118000000),
    @(45072, 7713658000000),
    @(45079, 7712715000000),
    @(45086, 7714391000000),
    @(45093, 7709739000000),
    @(45100, 7710607000000),
    @(45107, 7219695000000),
    @(45114, 7206902000000),
    @(45121, 7205494000000)
)

$lastRow = 1271
$startRow = $lastRow + 1

# Copy the number format/style of the last existing date cell (A1271)
# down onto the new date cells before filling in the values.
$formatSource = $ws.Cells.Item($lastRow, 1)
$formatSource.Copy()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $date = $data[$i][0]
    $val = $data[$i][1]

    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = ":ECBASSETSW"
    $ws.Cells.Item($r, 3).Value = $val
    $ws.Cells.Item($r, 4).Value = $val
    $ws.Cells.Item($r, 5).Value = $val
    $ws.Cells.Item($r, 6).Value = $val
    $ws.Cells.Item($r, 7).Value = 0
}
